{"js": "// Apply the three textual edits described by the diff:\n// 1. \"dox-ol\u00e1snak\" -> \"doxol\u00e1snak\"            (remove hyphen / spelling fix)\n// 2. \"legf\u0151bbk\u00e9ppen\" -> \"legf\u0151k\u00e9ppen\"          (remove duplicated \"b\" / spelling fix)\n// 3. Extend the final paragraph: replace the trailing\n//    \"\u00dagyhogy szerintem az internetes \" with a new closing passage.\n\nconst body = context.document.body;\n\n// --- Edit 1 -----------------------------------------------------------\nlet results1 = body.search(\"dox-ol\u00e1snak\", { matchCase: true });\nresults1.load(\"text\");\nawait context.sync();\nif (results1.items.length > 0) {\n  results1.items[0].insertText(\"doxol\u00e1snak\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Edit 2 -----------------------------------------------------------\nlet results2 = body.search(\"legf\u0151bbk\u00e9ppen\", { matchCase: true });\nresults2.load(\"text\");\nawait context.sync();\nif (results2.items.length > 0) {\n  results2.items[0].insertText(\"legf\u0151k\u00e9ppen\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Edit 3 -------------------------------------------------------------\nconst addition =\n  \"Term\u00e9szetesen ez csak egy sz\u00e9ls\u0151s\u00e9ges p\u00e9lda, ugyanis szem\u00e9lyis\u00e9gf\u00fcgg\u0151 az, \" +\n  \"hogy ki hogyan kezeli az ilyen helyzeteket, \u00e9s lehets\u00e9ges, hogy pont egy \" +\n  \"konfliktus miatt fog kialakulni valaki k\u00f6z\u00f6tt egy bar\u00e1ts\u00e1g az interneten, \" +\n  \"viszont az is kijelenthet\u0151, hogy ez az olyan platformokon, mint a Twitter, \" +\n  \"Facebook, Instagram sokkal kisebb es\u00e9ly van a felhaszn\u00e1l\u00f3b\u00e1zis  \" +\n  \"ellens\u00e9gess\u00e9ge miatt, melyekr\u0151l t\u00f6bb tanulm\u00e1ny is k\u00e9sz\u00fclt. \";\n\nlet results3 = body.search(\"\u00dagyhogy szerintem az internetes \", { matchCase: true });\nresults3.load(\"text\");\nawait context.sync();\nif (results3.items.length > 0) {\n  results3.items[0].insertText(addition, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply the three textual edits described by the diff:\n# 1. \"dox-ol\u00e1snak\" -> \"doxol\u00e1snak\"            (remove hyphen / spelling fix)\n# 2. \"legf\u0151bbk\u00e9ppen\" -> \"legf\u0151k\u00e9ppen\"          (remove duplicated \"b\" / spelling fix)\n# 3. Extend the final paragraph: replace the trailing\n#    \"\u00dagyhogy szerintem az internetes \" with a new closing passage.\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n# --- Edit 1 -------------------------------------------------------------\n$find1 = $d.Content.Find\n$find1.Execute(\n  \"dox-ol\u00e1snak\",   # FindText\n  $true,           # MatchCase\n  $false,          # MatchWholeWord\n  $false,          # MatchWildcards\n  $false,          # MatchSoundsLike\n  $false,          # MatchAllWordForms\n  $true,           # Forward\n  $wdFindContinue, # Wrap\n  $false,          # Format\n  \"doxol\u00e1snak\",    # ReplaceWith\n  $wdReplaceAll    # Replace\n)\n\n# --- Edit 2 -------------------------------------------------------------\n$find2 = $d.Content.Find\n$find2.Execute(\n  \"legf\u0151bbk\u00e9ppen\", # FindText\n  $true,           # MatchCase\n  $false,          # MatchWholeWord\n  $false,          # MatchWildcards\n  $false,          # MatchSoundsLike\n  $false,          # MatchAllWordForms\n  $true,           # Forward\n  $wdFindContinue, # Wrap\n  $false,          # Format\n  \"legf\u0151k\u00e9ppen\",   # ReplaceWith\n  $wdReplaceAll    # Replace\n)\n\n# --- Edit 3 -------------------------------------------------------------\n$addition = \"Term\u00e9szetesen ez csak egy sz\u00e9ls\u0151s\u00e9ges p\u00e9lda, ugyanis szem\u00e9lyis\u00e9gf\u00fcgg\u0151 az, hogy ki hogyan kezeli az ilyen helyzeteket, \u00e9s lehets\u00e9ges, hogy pont egy konfliktus miatt fog kialakulni valaki k\u00f6z\u00f6tt egy bar\u00e1ts\u00e1g az interneten, viszont az is kijelenthet\u0151, hogy ez az olyan platformokon, mint a Twitter, Facebook, Instagram sokkal kisebb es\u00e9ly van a felhaszn\u00e1l\u00f3b\u00e1zis  ellens\u00e9gess\u00e9ge miatt, melyekr\u0151l t\u00f6bb tanulm\u00e1ny is k\u00e9sz\u00fclt. \"\n\n$find3 = $d.Content.Find\n$find3.Execute(\n  \"\u00dagyhogy szerintem az internetes \", # FindText\n  $true,           # MatchCase\n  $false,          # MatchWholeWord\n  $false,          # MatchWildcards\n  $false,          # MatchSoundsLike\n  $false,          # MatchAllWordForms\n  $true,           # Forward\n  $wdFindContinue, # Wrap\n  $false,          # Format\n  $addition,       # ReplaceWith\n  $wdReplaceAll    # Replace\n)\n"}
